$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The workbook's Sheet1 row 4 (JobHistory test-suite row) pulls its counts from
# an external workbook (JobHistory/_Test_Suite_Statistics.xlsx, external
# reference [3]) via cached-value formulas:
#   C4 =[3]Sheet1!$G$2   (errors)
#   D4 =[3]Sheet1!$E$1   (total)
#   G4 =[3]Sheet1!$G$5   (automated/tagged)
#   H4 =[3]Sheet1!$G$4   (automateable)
# The source workbook now reports one more completed/tagged item, so the
# cached external values need to be refreshed: total 4->5, errors 2->3,
# automateable 22->24, automated/tagged 15->17.
#
# Updating these four source cells lets every downstream formula on the
# sheet (E4, I4, L1, N1, N2, L3, N3, P3, L5, L6, L7, ...) recompute itself
# automatically from its existing formula - only the raw inputs change.
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 5
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 24

$wb.RefreshAll()
$excel.CalculateFullRebuild()
